$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column N (MACRO_SCORE) for rows 2-6 with the new recalculated value
$ws.Range("N2:N6").Value = 85.77505782882612
